# Update ER and EDC flow: add a "tenant_id" field to the Collect_spending
# table, right before the existing "tenant_name" column.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Collect_spending")

# Make this the active sheet/tab, matching the workbook's tabSelected state.
[void]$ws.Activate()

# Insert a new column before column G (7th column), shifting
# tenant_name/receipt_id/... and everything to its right one column over.
[void]$ws.Columns.Item(7).Insert()

# Populate the new column's header (row 1) and type (row 2) cells.
$ws.Range("G1").Value = "tenant_id"
$ws.Range("G2").Value = "string"

# Leave the selection on the newly added header's data cell, as in the
# final workbook.
[void]$ws.Range("G2").Select()
